$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The account-statement table currently has a single data row (row 16, period
# "2507"). A new period ("2508") is being added for the same worker, so we
# insert a fresh row right below the existing data row (row 17) - this
# pushes the trailing "signature" block (rows 21-22) down to rows 22-23,
# exactly like the target workbook.
$ws.Rows("17").Insert()

# Duplicate row 16 (worker CC/15676416/ELKIN JAVIER CUARTAS TORRES, with its
# Salario Basico / Valor Mora figures) into the newly inserted row 17,
# copying values + formatting together.
$ws.Range("B16:J16").Copy($ws.Range("B17:J17"))

# The new row represents period 2508 instead of 2507.
$ws.Range("E17").Value = "2508"

# Update the "Cant. Periodos" counter: there are now 2 periods in arrears.
$ws.Range("F13").Value = 2

# Update the total "VALOR MORA" (sum of the Valor Mora column: 56940 + 56940).
$ws.Range("E11").Value = 113880
